# Swap the data (columns B and E:AB) between row pairs (25,26) and (85,86).
# Columns A (id), C (Div) and D (Date) stay attached to their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2)

    # Column B (numeric match id) plus columns E..AB swap between the two rows.
    $cols = @("B") + @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$row1")
        $cell2 = $ws.Range("$col$row2")
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-RowData $ws 25 26
Swap-RowData $ws 85 86
